$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 89.38217433333334
$ws.Range("H2").Value = 268.146523
$ws.Range("I2").Value = 0.2143552015363441
$ws.Range("J2").Value = 0.2175965347165783
$ws.Range("M2").Value = 2.906846333333333
$ws.Range("N2").Value = 8.720538999999999
$ws.Range("O2").Value = 0.005520525738044089
$ws.Range("P2").Value = 0.005624540846623205
$ws.Range("Q2").Value = 259.8202457262107
$ws.Range("R2").Value = 2338.382211535897
$ws.Range("S2").Value = 0.001183353407165015
$ws.Range("T2").Value = 0.001223880597597059
$ws.Range("G3").Value = 89.38217433333334
$ws.Range("H3").Value = 268.146523
$ws.Range("I3").Value = 0.2143552015363441
$ws.Range("J3").Value = 0.2175965347165783
$ws.Range("O3").Value = 0.3528665483720876
$ws.Range("P3").Value = 0.3595150912979765
$ws.Range("Q3").Value = 16607.45328561376
$ws.Range("R3").Value = 149467.0795705239
$ws.Range("S3").Value = 0.07563878009173296
$ws.Range("T3").Value = 0.07822923804475398
$ws.Range("G4").Value = 89.38217433333334
$ws.Range("H4").Value = 268.146523
$ws.Range("I4").Value = 0.2143552015363441
$ws.Range("J4").Value = 0.2175965347165783
$ws.Range("M4").Value = 137.0717086666666
$ws.Range("N4").Value = 411.2151259999999
$ws.Range("O4").Value = 0.2603191943704447
$ws.Range("P4").Value = 0.2652240042658267
$ws.Range("Q4").Value = 12251.76736021188
$ws.Range("R4").Value = 110265.9062419069
$ws.Range("S4").Value = 0.0558007733730554
$ws.Range("T4").Value = 0.05771182425189888
$ws.Range("G5").Value = 89.38217433333334
$ws.Range("H5").Value = 268.146523
$ws.Range("I5").Value = 0.2143552015363441
$ws.Range("J5").Value = 0.2175965347165783
$ws.Range("M5").Value = 29.2127365
$ws.Range("N5").Value = 58.425473
$ws.Range("O5").Value = 0.05547925319534149
$ws.Range("P5").Value = 0.03768304451958546
$ws.Range("Q5").Value = 2611.09790659673
$ws.Range("R5").Value = 15666.58743958038
$ws.Range("S5").Value = 0.01189226649977329
$ws.Range("T5").Value = 0.008199699905032344
$ws.Range("G6").Value = 89.38217433333334
$ws.Range("H6").Value = 268.146523
$ws.Range("I6").Value = 0.2143552015363441
$ws.Range("J6").Value = 0.2175965347165783
$ws.Range("M6").Value = 171.5584106666666
$ws.Range("N6").Value = 514.6752319999999
$ws.Range("O6").Value = 0.3258144783240821
$ws.Range("P6").Value = 0.331953319069988
$ws.Range("Q6").Value = 15334.26377055759
$ws.Range("R6").Value = 138008.3739350183
$ws.Range("S6").Value = 0.06984002816461744
$ws.Range("T6").Value = 0.07223189191729605
$ws.Range("I7").Value = 0.2934277926151677
$ws.Range("J7").Value = 0.2978648075949286
$ws.Range("M7").Value = 2.906846333333333
$ws.Range("N7").Value = 8.720538999999999
$ws.Range("O7").Value = 0.005520525738044089
$ws.Range("P7").Value = 0.005624540846623205
$ws.Range("Q7").Value = 355.6642462312545
$ws.Range("R7").Value = 3200.97821608129
$ws.Range("S7").Value = 0.001619875681389496
$ws.Range("T7").Value = 0.001675352777089238
$ws.Range("I8").Value = 0.2934277926151677
$ws.Range("J8").Value = 0.2978648075949286
$ws.Range("O8").Value = 0.3528665483720876
$ws.Range("P8").Value = 0.3595150912979765
$ws.Range("S8").Value = 0.103540852376555
$ws.Range("T8").Value = 0.107086893496945
$ws.Range("I9").Value = 0.2934277926151677
$ws.Range("J9").Value = 0.2978648075949286
$ws.Range("M9").Value = 137.0717086666666
$ws.Range("N9").Value = 411.2151259999999
$ws.Range("O9").Value = 0.2603191943704447
$ws.Range("P9").Value = 0.2652240042658267
$ws.Range("Q9").Value = 16771.27042579367
$ws.Range("R9").Value = 150941.4338321431
$ws.Range("S9").Value = 0.07638488657947835
$ws.Range("T9").Value = 0.07900089700019698
$ws.Range("I10").Value = 0.2934277926151677
$ws.Range("J10").Value = 0.2978648075949286
$ws.Range("M10").Value = 29.2127365
$ws.Range("N10").Value = 58.425473
$ws.Range("O10").Value = 0.05547925319534149
$ws.Range("P10").Value = 0.03768304451958546
$ws.Range("Q10").Value = 3574.294859856056
$ws.Range("R10").Value = 21445.76915913633
$ws.Range("S10").Value = 0.01627915480104704
$ws.Range("T10").Value = 0.01122445280541745
$ws.Range("I11").Value = 0.2934277926151677
$ws.Range("J11").Value = 0.2978648075949286
$ws.Range("M11").Value = 171.5584106666666
$ws.Range("N11").Value = 514.6752319999999
$ws.Range("O11").Value = 0.3258144783240821
$ws.Range("P11").Value = 0.331953319069988
$ws.Range("Q11").Value = 20990.85600593909
$ws.Range("R11").Value = 188917.7040534518
$ws.Range("S11").Value = 0.0956030231766978
$ws.Range("T11").Value = 0.09887721151527992
$ws.Range("G12").Value = 90.33462533333334
$ws.Range("H12").Value = 271.003876
$ws.Range("I12").Value = 0.2166393574945233
$ws.Range("J12").Value = 0.2199152301234996
$ws.Range("M12").Value = 2.906846333333333
$ws.Range("N12").Value = 8.720538999999999
$ws.Range("O12").Value = 0.005520525738044089
$ws.Range("P12").Value = 0.005624540846623205
$ws.Range("Q12").Value = 262.5888744232404
$ws.Range("R12").Value = 2363.299869809164
$ws.Range("S12").Value = 0.00119596314892185
$ws.Range("T12").Value = 0.001236922194624166
$ws.Range("G13").Value = 90.33462533333334
$ws.Range("H13").Value = 271.003876
$ws.Range("I13").Value = 0.2166393574945233
$ws.Range("J13").Value = 0.2199152301234996
$ws.Range("O13").Value = 0.3528665483720876
$ws.Range("P13").Value = 0.3595150912979765
$ws.Range("Q13").Value = 16784.42129525679
$ws.Range("R13").Value = 151059.7916573111
$ws.Range("S13").Value = 0.07644478232063918
$ws.Range("T13").Value = 0.07906284403566549
$ws.Range("G14").Value = 90.33462533333334
$ws.Range("H14").Value = 271.003876
$ws.Range("I14").Value = 0.2166393574945233
$ws.Range("J14").Value = 0.2199152301234996
$ws.Range("M14").Value = 137.0717086666666
$ws.Range("N14").Value = 411.2151259999999
$ws.Range("O14").Value = 0.2603191943704447
$ws.Range("P14").Value = 0.2652240042658267
$ws.Range("Q14").Value = 12382.32144620315
$ws.Range("R14").Value = 111440.8930158284
$ws.Range("S14").Value = 0.05639538301190505
$ws.Range("T14").Value = 0.05832679793239533
$ws.Range("G15").Value = 90.33462533333334
$ws.Range("H15").Value = 271.003876
$ws.Range("I15").Value = 0.2166393574945233
$ws.Range("J15").Value = 0.2199152301234996
$ws.Range("M15").Value = 29.2127365
$ws.Range("N15").Value = 58.425473
$ws.Range("O15").Value = 0.05547925319534149
$ws.Range("P15").Value = 0.03768304451958546
$ws.Range("Q15").Value = 2638.921606688891
$ws.Range("R15").Value = 15833.52964013335
$ws.Range("S15").Value = 0.01201898976651476
$ws.Range("T15").Value = 0.008287075407278717
$ws.Range("G16").Value = 90.33462533333334
$ws.Range("H16").Value = 271.003876
$ws.Range("I16").Value = 0.2166393574945233
$ws.Range("J16").Value = 0.2199152301234996
$ws.Range("M16").Value = 171.5584106666666
$ws.Range("N16").Value = 514.6752319999999
$ws.Range("O16").Value = 0.3258144783240821
$ws.Range("P16").Value = 0.331953319069988
$ws.Range("Q16").Value = 15497.66475035547
$ws.Range("R16").Value = 139478.9827531992
$ws.Range("S16").Value = 0.07058423924654242
$ws.Range("T16").Value = 0.07300159055353592
$ws.Range("G17").Value = 18.634161
$ws.Range("H17").Value = 37.268322
$ws.Range("I17").Value = 0.0446882095496985
$ws.Range("J17").Value = 0.03024263611988591
$ws.Range("M17").Value = 2.906846333333333
$ws.Range("N17").Value = 8.720538999999999
$ws.Range("O17").Value = 0.005520525738044089
$ws.Range("P17").Value = 0.005624540846623205
$ws.Range("Q17").Value = 54.16664257759299
$ws.Range("R17").Value = 324.9998554655579
$ws.Range("S17").Value = 0.0002467024110062182
$ws.Range("T17").Value = 0.0001701009421658606
$ws.Range("G18").Value = 18.634161
$ws.Range("H18").Value = 37.268322
$ws.Range("I18").Value = 0.0446882095496985
$ws.Range("J18").Value = 0.03024263611988591
$ws.Range("O18").Value = 0.3528665483720876
$ws.Range("P18").Value = 0.3595150912979765
$ws.Range("Q18").Value = 3462.27825214917
$ws.Range("R18").Value = 20773.66951289502
$ws.Range("S18").Value = 0.01576897425673068
$ws.Range("T18").Value = 0.01087268408573227
$ws.Range("G19").Value = 18.634161
$ws.Range("H19").Value = 37.268322
$ws.Range("I19").Value = 0.0446882095496985
$ws.Range("J19").Value = 0.03024263611988591
$ws.Range("M19").Value = 137.0717086666666
$ws.Range("N19").Value = 411.2151259999999
$ws.Range("O19").Value = 0.2603191943704447
$ws.Range("P19").Value = 0.2652240042658267
$ws.Range("Q19").Value = 2554.216287839761
$ws.Range("R19").Value = 15325.29772703857
$ws.Range("S19").Value = 0.01163319870783512
$ws.Range("T19").Value = 0.008021073051270467
$ws.Range("G20").Value = 18.634161
$ws.Range("H20").Value = 37.268322
$ws.Range("I20").Value = 0.0446882095496985
$ws.Range("J20").Value = 0.03024263611988591
$ws.Range("M20").Value = 29.2127365
$ws.Range("N20").Value = 58.425473
$ws.Range("O20").Value = 0.05547925319534149
$ws.Range("P20").Value = 0.03768304451958546
$ws.Range("Q20").Value = 544.3548351915764
$ws.Range("R20").Value = 2177.419340766306
$ws.Range("S20").Value = 0.0024792684924542
$ws.Range("T20").Value = 0.001139634603295284
$ws.Range("G21").Value = 18.634161
$ws.Range("H21").Value = 37.268322
$ws.Range("I21").Value = 0.0446882095496985
$ws.Range("J21").Value = 0.03024263611988591
$ws.Range("M21").Value = 171.5584106666666
$ws.Range("N21").Value = 514.6752319999999
$ws.Range("O21").Value = 0.3258144783240821
$ws.Range("P21").Value = 0.331953319069988
$ws.Range("Q21").Value = 3196.847045266783
$ws.Range("R21").Value = 19181.0822716007
$ws.Range("S21").Value = 0.01456006568167228
$ws.Range("T21").Value = 0.01003914343742203
$ws.Range("G22").Value = 96.27664699999998
$ws.Range("H22").Value = 288.829941
$ws.Range("I22").Value = 0.2308894388042666
$ws.Range("J22").Value = 0.2343807914451077
$ws.Range("M22").Value = 2.906846333333333
$ws.Range("N22").Value = 8.720538999999999
$ws.Range("O22").Value = 0.005520525738044089
$ws.Range("P22").Value = 0.005624540846623205
$ws.Range("Q22").Value = 279.8614183175775
$ws.Range("R22").Value = 2518.752764858198
$ws.Range("S22").Value = 0.001274631089561509
$ws.Range("T22").Value = 0.001318284335146883
$ws.Range("G23").Value = 96.27664699999998
$ws.Range("H23").Value = 288.829941
$ws.Range("I23").Value = 0.2308894388042666
$ws.Range("J23").Value = 0.2343807914451077
$ws.Range("O23").Value = 0.3528665483720876
$ws.Range("P23").Value = 0.3595150912979765
$ws.Range("Q23").Value = 17888.46522781158
$ws.Range("R23").Value = 160996.1870503042
$ws.Range("S23").Value = 0.0814731593264299
$ws.Range("T23").Value = 0.08426343163487987
$ws.Range("G24").Value = 96.27664699999998
$ws.Range("H24").Value = 288.829941
$ws.Range("I24").Value = 0.2308894388042666
$ws.Range("J24").Value = 0.2343807914451077
$ws.Range("M24").Value = 137.0717086666666
$ws.Range("N24").Value = 411.2151259999999
$ws.Range("O24").Value = 0.2603191943704447
$ws.Range("P24").Value = 0.2652240042658267
$ws.Range("Q24").Value = 13196.8045089875
$ws.Range("R24").Value = 118771.2405808875
$ws.Range("S24").Value = 0.06010495269817075
$ws.Range("T24").Value = 0.06216341203006507
$ws.Range("G25").Value = 96.27664699999998
$ws.Range("H25").Value = 288.829941
$ws.Range("I25").Value = 0.2308894388042666
$ws.Range("J25").Value = 0.2343807914451077
$ws.Range("M25").Value = 29.2127365
$ws.Range("N25").Value = 58.425473
$ws.Range("O25").Value = 0.05547925319534149
$ws.Range("P25").Value = 0.03768304451958546
$ws.Range("Q25").Value = 2812.504319914515
$ws.Range("R25").Value = 16875.02591948709
$ws.Range("S25").Value = 0.01280957363555221
$ws.Range("T25").Value = 0.008832181798561666
$ws.Range("G26").Value = 96.27664699999998
$ws.Range("H26").Value = 288.829941
$ws.Range("I26").Value = 0.2308894388042666
$ws.Range("J26").Value = 0.2343807914451077
$ws.Range("M26").Value = 171.5584106666666
$ws.Range("N26").Value = 514.6752319999999
$ws.Range("O26").Value = 0.3258144783240821
$ws.Range("P26").Value = 0.331953319069988
$ws.Range("Q26").Value = 16517.0685436357
$ws.Range("R26").Value = 148653.6168927213
$ws.Range("S26").Value = 0.07522712205455219
$ws.Range("T26").Value = 0.07780348164645415
